# Updates cryptos list values (Price column D and Volume(1h) column E)
# per the commit "Updated cryptos list on Mon Aug 12 14:15:31 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.691.38'
$ws.Range('E2').Value = '  -2.76%  '
$ws.Range('D3').Value = '2.620.45'
$ws.Range('E3').Value = '  -0.39%  '
$ws.Range('D4').Value = "'0.997"
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.27%  '
$ws.Range('D5').Value = "'510.07"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.63%  '
$ws.Range('D6').Value = "'145.68"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -3.67%  '
$ws.Range('D7').Value = "'0.994"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.45%  '
$ws.Range('D8').Value = "'0.566"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.35%  '
$ws.Range('D9').Value = '2.638.75'
$ws.Range('E9').Value = '  +0.12%  '
$ws.Range('D10').Value = "'6.49"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +1.56%  '
$ws.Range('D11').Value = "'0.104"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -2.12%  '
$ws.Range('D12').Value = "'0.335"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -2.57%  '
$ws.Range('E13').Value = '  -1.48%  '
$ws.Range('D14').Value = '3.079.94'
$ws.Range('E14').Value = '  -0.33%  '
$ws.Range('D15').Value = '58.536.64'
$ws.Range('E15').Value = '  -3.04%  '
$ws.Range('D16').Value = "'21.01"
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -2.33%  '
$ws.Range('D17').Value = "'0.0000137"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.45%  '
$ws.Range('D18').Value = '2.634.59'
$ws.Range('E18').Value = '  -0.08%  '
$ws.Range('D19').Value = "'4.54"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.41%  '
$ws.Range('D20').Value = "'341.13"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.63%  '
$ws.Range('D21').Value = "'10.33"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.23%  '
$ws.Range('D22').Value = "'6.09"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.39%  '
$ws.Range('D23').Value = "'0.997"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.25%  '
$ws.Range('D24').Value = "'60.28"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.07%  '
$ws.Range('D25').Value = "'0.421"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.20%  '
$ws.Range('D26').Value = '2.726.11'
$ws.Range('E26').Value = '  -0.67%  '
$ws.Range('D27').Value = "'0.995"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.56%  '
$ws.Range('D28').Value = "'0.159"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.78%  '
$ws.Range('D29').Value = '0.0₃0814'
$ws.Range('E29').Value = '  -1.76%  '
$ws.Range('D30').Value = "'7.06"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.11%  '
$ws.Range('D31').Value = "'0.997"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.32%  '
$ws.Range('D32').Value = "'6.53"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +7.97%  '
$ws.Range('D33').Value = "'18.81"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.26%  '
$ws.Range('D34').Value = "'1.56"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -2.23%  '
$ws.Range('D35').Value = "'148.52"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.38%  '
$ws.Range('D36').Value = "'1.04"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +15.86%  '
$ws.Range('D37').Value = "'3.97"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.74%  '
$ws.Range('D38').Value = "'1.14"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.16%  '
$ws.Range('D39').Value = "'0.855"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.82%  '
$ws.Range('D40').Value = "'36.23"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.21%  '
$ws.Range('D41').Value = "'3.66"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.25%  '
$ws.Range('D42').Value = "'1.41"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -3.29%  '
$ws.Range('B43').Value = 'Bittensor'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D43').Value = "'281.41"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -3.17%  '
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').Value = "'0.616"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.48%  '
$ws.Range('B45').Value = 'Stellar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D45').Value = "'0.0988"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.76%  '
$ws.Range('B46').Value = 'FirstDigitalUSD'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D46').Value = "'0.991"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.53%  '
$ws.Range('D47').Value = "'19.51"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.35%  '
$ws.Range('D48').Value = "'0.0537"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.77%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').Value = "'0.0230"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.27%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').Value = "'4.72"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.38%  '
$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').Value = "'10.24"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.47%  '
